$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new log entry as row 59, mirroring the columns used by the
# existing rows: B=Start, C=Stop, E=Delta, F=Activity Category, G=Activity Summary
$ws.Range("B59").Value = "1:30PM"
$ws.Range("C59").Value = "1:42PM"
$ws.Range("E59").Value = 12
$ws.Range("F59").Value = "Debug"
$ws.Range("G59").Value = "Debugging the New Order GUI response to changing product amounts while in-order list"

# Selection state: select the new last cell in the summary column, matching
# the view being scrolled down to the newly added row.
$ws.Range("G59").Select()
